# This script applies the odds updates described in the commit diff
# (FlashScore 2025-05-19 weekly games odds refresh).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Argentinos Jrs vs San Lorenzo
$ws.Range("G3").Value = 1.75
$ws.Range("H3").Value = 3
$ws.Range("I3").Value = 6.5
$ws.Range("U3").Value = 6.5
$ws.Range("V3").Value = 10
$ws.Range("W3").Value = 13
$ws.Range("AA3").Value = 6.5
$ws.Range("AB3").Value = 26
$ws.Range("AC3").Value = 126
$ws.Range("AE3").Value = 10
$ws.Range("AF3").Value = 29
$ws.Range("AG3").Value = 23
$ws.Range("AI3").Value = 67
$ws.Range("AJ3").Value = 81

# Row 4: Boca Juniors vs Independiente
$ws.Range("H4").Value = 2.9
$ws.Range("I4").Value = 3.25
$ws.Range("L4").Value = 1.53
$ws.Range("M4").Value = 2.38
$ws.Range("N4").Value = 2.7
$ws.Range("O4").Value = 1.44
$ws.Range("P4").Value = 1.57
$ws.Range("Q4").Value = 2.25
$ws.Range("AA4").Value = 5.5
$ws.Range("AF4").Value = 15

# Row 6: Agropecuario vs Gimnasia Jujuy
$ws.Range("G6").Value = 2.7
$ws.Range("H6").Value = 2.75
$ws.Range("I6").Value = 2.9
$ws.Range("T6").Value = 5.5
$ws.Range("Z6").Value = 5
$ws.Range("AF6").Value = 12
$ws.Range("AG6").Value = 13
$ws.Range("AH6").Value = 34

# Row 10: Colo Colo vs Nublense
$ws.Range("L10").Value = 1.25
$ws.Range("M10").Value = 3.75
$ws.Range("N10").Value = 1.9
$ws.Range("O10").Value = 1.9

# Row 11: La Serena vs Deportes Iquique
$ws.Range("K11").Value = 8

# Row 12: Llaneros vs Deportes Tolima
$ws.Range("G12").Value = 2.9
$ws.Range("I12").Value = 2.45
$ws.Range("J12").Value = 1.11
$ws.Range("K12").Value = 6.5
$ws.Range("N12").Value = 2.6
$ws.Range("O12").Value = 1.48
$ws.Range("P12").Value = 1.53
$ws.Range("Q12").Value = 2.38
$ws.Range("T12").Value = 7
$ws.Range("U12").Value = 13
$ws.Range("V12").Value = 12
$ws.Range("W12").Value = 34
$ws.Range("X12").Value = 29
$ws.Range("AD12").Value = 1250
$ws.Range("AE12").Value = 6.5
$ws.Range("AF12").Value = 11
$ws.Range("AH12").Value = 23
$ws.Range("AI12").Value = 26

# Row 15: LDU Quito vs Libertad
$ws.Range("N15").Value = 1.7
$ws.Range("O15").Value = 2.1

# Row 26: Binacional vs Cienciano
$ws.Range("G26").Value = 3.5
$ws.Range("H26").Value = 3.4
$ws.Range("I26").Value = 2.05
$ws.Range("J26").Value = 1.04
$ws.Range("K26").Value = 13
$ws.Range("R26").Value = 1.67
$ws.Range("S26").Value = 2.1
$ws.Range("T26").Value = 12
$ws.Range("U26").Value = 19
$ws.Range("V26").Value = 12
$ws.Range("W26").Value = 41
$ws.Range("X26").Value = 26
$ws.Range("Y26").Value = 34
$ws.Range("Z26").Value = 12
$ws.Range("AB26").Value = 13
$ws.Range("AE26").Value = 9
$ws.Range("AF26").Value = 11
$ws.Range("AG26").Value = 9
$ws.Range("AH26").Value = 19
$ws.Range("AI26").Value = 15
$ws.Range("AJ26").Value = 23

# Row 27: Alianza Lima vs Alianza Huanuco
$ws.Range("G27").Value = 1.22
$ws.Range("H27").Value = 5.5
$ws.Range("I27").Value = 13
$ws.Range("L27").Value = 1.17
$ws.Range("M27").Value = 5
$ws.Range("N27").Value = 1.57
$ws.Range("O27").Value = 2.35
$ws.Range("R27").Value = 2.25
$ws.Range("S27").Value = 1.57
$ws.Range("AA27").Value = 11
$ws.Range("AB27").Value = 26
$ws.Range("AF27").Value = 51
$ws.Range("AG27").Value = 34

# Row 47: Plaza Colonia vs Danubio
$ws.Range("G47").Value = 2.3
$ws.Range("H47").Value = 3
$ws.Range("I47").Value = 3
$ws.Range("L47").Value = 1.62
$ws.Range("M47").Value = 2.2
$ws.Range("P47").Value = 1.62
$ws.Range("Q47").Value = 2.2
$ws.Range("V47").Value = 11
$ws.Range("W47").Value = 23
$ws.Range("X47").Value = 26
$ws.Range("Z47").Value = 6
$ws.Range("AE47").Value = 6.5
$ws.Range("AF47").Value = 13

# Row 48: Boston River vs CA Cerro
$ws.Range("C48").Value = "20:15"
$ws.Range("J48").Value = 1.1
$ws.Range("K48").Value = 7
